$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert 4 new rows (for RULE-107..RULE-110) right after the two
#    existing PASS rows (RULE-105, RULE-106) and before the FAIL rows.
#    New rows inherit formatting from the row above (style "2" / PASS look).
# ------------------------------------------------------------------
$ws.Rows("4:7").Insert()

# RULE-107
$ws.Range("A4").Value = "RULE-107"
$ws.Range("B4").Value = "Forbidden substring check for .properties files"
$ws.Range("C4").Value = "HIGH"
$ws.Range("D4").Value = "PASS"
$ws.Range("E4").Value = "All checks passed"

# RULE-108
$ws.Range("A5").Value = "RULE-108"
$ws.Range("B5").Value = "Forbidden substring check for .policy files"
$ws.Range("C5").Value = "HIGH"
$ws.Range("D5").Value = "PASS"
$ws.Range("E5").Value = "All checks passed"

# RULE-109
$ws.Range("A6").Value = "RULE-109"
$ws.Range("B6").Value = "Forbidden regex pattern (ip addresses) check in .properties files"
$ws.Range("C6").Value = "HIGH"
$ws.Range("D6").Value = "PASS"
$ws.Range("E6").Value = "All checks passed"

# RULE-110
$ws.Range("A7").Value = "RULE-110"
$ws.Range("B7").Value = "Forbidden regex pattern (ip addresses) check in .policy files"
$ws.Range("C7").Value = "HIGH"
$ws.Range("D7").Value = "PASS"
$ws.Range("E7").Value = "All checks passed"

# ------------------------------------------------------------------
# 2. Update the wording of existing failure detail cells.
#    Row 8  = RULE-100 (was row 4): "Token" -> "Required token"
#    Row 9  = RULE-101 (was row 5): "[true, false, test]" -> "[true, false]"
#    Row 10 = RULE-102 (was row 6): "Token" -> "Required token"
# ------------------------------------------------------------------
$rule100Details = @"
• Validation failures:
• Required token 'apiId' not found in file: Properties\OCP\01\ITE.properties (case-sensitive: true)
• Required token 'apiId' not found in file: Properties\OCP\01\PROD.properties (case-sensitive: true)
• Required token 'apiId' not found in file: Properties\OCP\01\TDV.properties (case-sensitive: true)
• Required token 'apiId' not found in file: Properties\OCP\ITE.properties (case-sensitive: true)
• Required token 'apiId' not found in file: Properties\OCP\PROD.properties (case-sensitive: true)
• Required token 'apiId' not found in file: Properties\OCP\TDV.properties (case-sensitive: true)
"@

$rule101Details = @"
• Validation failures:
• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\01\ITE.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\01\ITE.properties
• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\01\PROD.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\01\PROD.properties
• Property 'LogJsonFormat' not found in file: Properties\OCP\01\TDV.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\01\TDV.properties
• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\ITE.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\ITE.properties
• Property 'LogJsonFormat' found but value does not match expected values [true, false] in file: Properties\OCP\PROD.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\PROD.properties
• Property 'LogJsonFormat' not found in file: Properties\OCP\TDV.properties
• Property 'anotherpropertycheck' not found in file: Properties\OCP\TDV.properties
"@

$rule102Details = @"
• Validation failures:
• Required token 'http.protocols=HTTPS' not found in file: Policies\TDV.policy (case-sensitive: true)
• Required token 'http.private.port=8081' not found in file: Policies\TDV.policy (case-sensitive: true)
"@

$ws.Range("E8").Value = $rule100Details.TrimEnd()
$ws.Range("E9").Value = $rule101Details.TrimEnd()
$ws.Range("E10").Value = $rule102Details.TrimEnd()

# ------------------------------------------------------------------
# 3. Widen column B to fit the longer rule-name text (42.18 -> 52.0 px).
#    51.15 characters of COM ColumnWidth serializes to exactly width="52".
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 51.15
